$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Helper: renaming an InlineShape living in a Footer story can throw
# "addressed block not found" when the shape is fetched straight off
# Footer.Range.InlineShapes (a stale-handle quirk once the footer has
# more than one paragraph). Selecting the shape first and then reaching
# it through $word.Selection.InlineShapes works around that reliably,
# and is harmless (a no-op selection change) for headers too.
function Rename-InlineLogo($inlineShape, $newName) {
    $inlineShape.Select()
    $selected = $word.Selection.InlineShapes.Item(1)
    $selected.Name = $newName
}

# Headers: BTec logo picture, currently named "image2.jpg" -> "image1.jpg"
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$headerPrimary = $sec.Headers.Item(1)
Rename-InlineLogo $headerPrimary.Range.InlineShapes.Item(1) "image1.jpg"

$headerFirstPage = $sec.Headers.Item(2)
Rename-InlineLogo $headerFirstPage.Range.InlineShapes.Item(1) "image1.jpg"

# Footers: Pearson Edexcel logo picture, currently named "image1.png" -> "image2.png"
$footerPrimary = $sec.Footers.Item(1)
Rename-InlineLogo $footerPrimary.Range.InlineShapes.Item(1) "image2.png"

$footerFirstPage = $sec.Footers.Item(2)
Rename-InlineLogo $footerFirstPage.Range.InlineShapes.Item(1) "image2.png"
